# Commit: "+ add calibration screeen"
#
# This adds a new "calibration screen" set of translation text entries to the
# Translation sheet (Text ID / Typography Name / Alignment / Direction / IT),
# and records the Wildcard Ranges value used by the "Default" typography on
# the Typography sheet.

$wb = $excel.ActiveWorkbook

$typo = $wb.Worksheets.Item("Typography")
$trans = $wb.Worksheets.Item("Translation")

# --- Typography sheet --------------------------------------------------
# The "Default" typography (row 4) now declares the wildcard range used by
# the new numeric/alphanumeric calibration widgets.
$typo.Range("I4").Value = "0-9,A-Z"

# --- Translation sheet --------------------------------------------------
# The language column header changes from "GB" to "IT" (Italian).
$trans.Range("F3").Value = "IT"

# New row: SingleUseId2 (generic placeholder value text)
$trans.Range("B5").Value = "SingleUseId2"
$trans.Range("C5").Value = "Default"
$trans.Range("D5").Value = "Center"
$trans.Range("E5").Value = "LTR"
$trans.Range("F5").Value = "<value>"

# New calibration screen texts
$trans.Range("B6").Value = "CALIB_PRESS_TOPSX"
$trans.Range("C6").Value = "Default"
$trans.Range("D6").Value = "Left"
$trans.Range("E6").Value = "LTR"
$trans.Range("F6").Value = "Toccare il crocino in alto per calibrare"

$trans.Range("B7").Value = "CALIB_PRESS_MIDDLEDX"
$trans.Range("C7").Value = "Default"
$trans.Range("D7").Value = "Left"
$trans.Range("E7").Value = "LTR"
$trans.Range("F7").Value = "Toccare il crocino al lato per calibrare"

$trans.Range("B8").Value = "CALIB_PRESS_BOTTOMCR"
$trans.Range("C8").Value = "Default"
$trans.Range("D8").Value = "Left"
$trans.Range("E8").Value = "LTR"
$trans.Range("F8").Value = "Toccare il crocino sotto per calibrare"

$trans.Range("B9").Value = "CALIB_DONE_FAIL"
$trans.Range("C9").Value = "Default"
$trans.Range("D9").Value = "Left"
$trans.Range("E9").Value = "LTR"
$trans.Range("F9").Value = "Calibrazione non riuscita premi per uscire"

$trans.Range("B10").Value = "CALIB_DONE_SUCCESS"
$trans.Range("C10").Value = "Default"
$trans.Range("D10").Value = "Left"
$trans.Range("E10").Value = "LTR"
$trans.Range("F10").Value = "Premere per salvare ed uscire"

$trans.Range("B11").Value = "CALIB_START_RELEASE"
$trans.Range("C11").Value = "Default"
$trans.Range("D11").Value = "Left"
$trans.Range("E11").Value = "LTR"
$trans.Range("F11").Value = "Rilascia per iniziare"

# Leave the cursor where the author last left it while entering this data.
[void]$trans.Activate()
[void]$trans.Range("C16").Select()
